$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1036.5
$ws.Range("J17").Value = 1061.3137
$ws.Range("L17").Value = 3183.9411
$ws.Range("N17").Value = -3519.9411
$ws.Range("H19").Value = 4055.8823
$ws.Range("I19").Value = 5495
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 5495
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -5320
$ws.Range("N19").Value = -2350
$ws.Range("H100").Value = 4687.4736
$ws.Range("I100").Value = 2771.2727
$ws.Range("J100").Value = 5468.148
$ws.Range("K100").Value = 2771.2727
$ws.Range("L100").Value = 5468.148
$ws.Range("M100").Value = -2230.2727
$ws.Range("N100").Value = -6550.148
$ws.Range("H127").Value = 699.8333
$ws.Range("I127").Value = 499.66666
$ws.Range("J127").Value = 900
$ws.Range("K127").Value = 1498.99998
$ws.Range("L127").Value = 2700
$ws.Range("M127").Value = 3461.00002
$ws.Range("N127").Value = -12620
$ws.Range("H132").Value = 2806.1538
$ws.Range("I132").Value = 1731.8572
$ws.Range("J132").Value = 5017.9414
$ws.Range("K132").Value = 5195.571599999999
$ws.Range("L132").Value = 15053.8242
$ws.Range("M132").Value = -2665.571599999999
$ws.Range("N132").Value = -20113.8242
$ws.Range("H137").Value = 2753.561
$ws.Range("I137").Value = 2692.9355
$ws.Range("K137").Value = 8078.806500000001
$ws.Range("M137").Value = -5528.806500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5119.3813
$ws.Range("I32").Value = 3042.9663
$ws.Range("J32").Value = 28219.5
$ws.Range("K32").Value = 3042.9663
$ws.Range("L32").Value = 28219.5
$ws.Range("M32").Value = -2755.9663
$ws.Range("N32").Value = -28793.5
$ws.Range("H61").Value = 440161.8
$ws.Range("I61").Value = 325286.03
$ws.Range("J61").Value = 719145.9
$ws.Range("K61").Value = 325286.03
$ws.Range("L61").Value = 719145.9
$ws.Range("M61").Value = -325074.03
$ws.Range("N61").Value = -719569.9
$ws.Range("H74").Value = 130970.05
$ws.Range("I74").Value = 152368.27
$ws.Range("J74").Value = 56639.367
$ws.Range("K74").Value = 152368.27
$ws.Range("L74").Value = 56639.367
$ws.Range("M74").Value = -151494.27
$ws.Range("N74").Value = -58387.367
$ws.Range("H76").Value = 25400
$ws.Range("J76").Value = 25400
$ws.Range("L76").Value = 25400
$ws.Range("N76").Value = -26076
$ws.Range("H77").Value = 130970.05
$ws.Range("I77").Value = 152368.27
$ws.Range("J77").Value = 56639.367
$ws.Range("K77").Value = 761841.35
$ws.Range("L77").Value = 283196.835
$ws.Range("M77").Value = -757473.35
$ws.Range("N77").Value = -291932.835
$ws.Range("H79").Value = 25400
$ws.Range("J79").Value = 25400
$ws.Range("L79").Value = 25400
$ws.Range("N79").Value = -27740
$ws.Range("H132").Value = 2732.9827
$ws.Range("I132").Value = 2483.1396
$ws.Range("J132").Value = 3449.2
$ws.Range("K132").Value = 7449.418799999999
$ws.Range("L132").Value = 10347.6
$ws.Range("M132").Value = -4919.418799999999
$ws.Range("N132").Value = -15407.6
$ws.Range("H136").Value = 440161.8
$ws.Range("I136").Value = 325286.03
$ws.Range("J136").Value = 719145.9
$ws.Range("K136").Value = 975858.0900000001
$ws.Range("L136").Value = 2157437.7
$ws.Range("M136").Value = -973308.0900000001
$ws.Range("N136").Value = -2162537.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1596.742
$ws.Range("I20").Value = 1346.909
$ws.Range("J20").Value = 2207.4443
$ws.Range("K20").Value = 1346.909
$ws.Range("L20").Value = 2207.4443
$ws.Range("M20").Value = -1099.909
$ws.Range("N20").Value = -2701.4443
$ws.Range("H86").Value = 5733.0713
$ws.Range("I86").Value = 8843.714
$ws.Range("J86").Value = 2622.4285
$ws.Range("K86").Value = 8843.714
$ws.Range("L86").Value = 2622.4285
$ws.Range("M86").Value = -7720.714
$ws.Range("N86").Value = -4868.4285
$ws.Range("H89").Value = 5733.0713
$ws.Range("I89").Value = 8843.714
$ws.Range("J89").Value = 2622.4285
$ws.Range("K89").Value = 44218.57
$ws.Range("L89").Value = 13112.1425
$ws.Range("M89").Value = -38602.57
$ws.Range("N89").Value = -24344.1425
$ws.Range("H105").Value = 1900.3636
$ws.Range("I105").Value = 1438.0769
$ws.Range("J105").Value = 3617.4285
$ws.Range("K105").Value = 1438.0769
$ws.Range("L105").Value = 3617.4285
$ws.Range("M105").Value = 308.9231
$ws.Range("N105").Value = -7111.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3002.6956
$ws.Range("I31").Value = 2000.6857
$ws.Range("J31").Value = 6190.909
$ws.Range("K31").Value = 2000.6857
$ws.Range("L31").Value = 6190.909
$ws.Range("M31").Value = -1705.6857
$ws.Range("N31").Value = -6780.909
$ws.Range("H34").Value = 3002.6956
$ws.Range("I34").Value = 2000.6857
$ws.Range("J34").Value = 6190.909
$ws.Range("K34").Value = 2000.6857
$ws.Range("L34").Value = 6190.909
$ws.Range("M34").Value = -1798.6857
$ws.Range("N34").Value = -6594.909
$ws.Range("H132").Value = 1893.7742
$ws.Range("I132").Value = 1042.3334
$ws.Range("J132").Value = 3681.8
$ws.Range("K132").Value = 3127.0002
$ws.Range("L132").Value = 11045.4
$ws.Range("M132").Value = -597.0001999999999
$ws.Range("N132").Value = -16105.4
$ws.Range("H134").Value = 1795.7906
$ws.Range("I134").Value = 1182.6666
$ws.Range("J134").Value = 2570.2632
$ws.Range("K134").Value = 3547.9998
$ws.Range("L134").Value = 7710.7896
$ws.Range("M134").Value = -1012.9998
$ws.Range("N134").Value = -12780.7896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5550
$ws.Range("I80").Value = 7203.5713
$ws.Range("J80").Value = 3445.4546
$ws.Range("K80").Value = 7203.5713
$ws.Range("L80").Value = 3445.4546
$ws.Range("M80").Value = -6205.5713
$ws.Range("N80").Value = -5441.4546
$ws.Range("H83").Value = 5550
$ws.Range("I83").Value = 7203.5713
$ws.Range("J83").Value = 3445.4546
$ws.Range("K83").Value = 36017.85649999999
$ws.Range("L83").Value = 17227.273
$ws.Range("M83").Value = -31025.85649999999
$ws.Range("N83").Value = -27211.273
$ws.Range("H132").Value = 4100.875
$ws.Range("I132").Value = 4164.3335
$ws.Range("J132").Value = 3979.7273
$ws.Range("K132").Value = 12493.0005
$ws.Range("L132").Value = 11939.1819
$ws.Range("M132").Value = -9963.000499999998
$ws.Range("N132").Value = -16999.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 777.94116
$ws.Range("I46").Value = 748.06665
$ws.Range("J46").Value = 1002
$ws.Range("K46").Value = 748.06665
$ws.Range("L46").Value = 1002
$ws.Range("M46").Value = -560.06665
$ws.Range("N46").Value = -1378
$ws.Range("H82").Value = 3065.3572
$ws.Range("I82").Value = 2414.1667
$ws.Range("J82").Value = 3553.75
$ws.Range("K82").Value = 2414.1667
$ws.Range("L82").Value = 3553.75
$ws.Range("M82").Value = -2053.1667
$ws.Range("N82").Value = -4275.75
$ws.Range("H85").Value = 3065.3572
$ws.Range("I85").Value = 2414.1667
$ws.Range("J85").Value = 3553.75
$ws.Range("K85").Value = 2414.1667
$ws.Range("L85").Value = 3553.75
$ws.Range("M85").Value = -1166.1667
$ws.Range("N85").Value = -6049.75
$ws.Range("H122").Value = 2349.5715
$ws.Range("I122").Value = 2299.4
$ws.Range("K122").Value = 6898.200000000001
$ws.Range("M122").Value = -4448.200000000001
$ws.Range("H132").Value = 6314.098
$ws.Range("I132").Value = 2095.6843
$ws.Range("J132").Value = 18644.846
$ws.Range("K132").Value = 6287.0529
$ws.Range("L132").Value = 55934.538
$ws.Range("M132").Value = -3757.0529
$ws.Range("N132").Value = -60994.538
$ws.Range("H136").Value = 3213.182
$ws.Range("I136").Value = 1722.6875
$ws.Range("K136").Value = 5168.0625
$ws.Range("M136").Value = -2618.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1179.8182
$ws.Range("I122").Value = 1047.8
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3143.4
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -693.3999999999996
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 1545.4529
$ws.Range("I132").Value = 1015.7838
$ws.Range("J132").Value = 2770.3125
$ws.Range("K132").Value = 3047.3514
$ws.Range("L132").Value = 8310.9375
$ws.Range("M132").Value = -517.3514
$ws.Range("N132").Value = -13370.9375
$ws.Range("H136").Value = 10317673
$ws.Range("I136").Value = 13903790
$ws.Range("J136").Value = 386888.5
$ws.Range("K136").Value = 41711370
$ws.Range("L136").Value = 1160665.5
$ws.Range("M136").Value = -41708820
$ws.Range("N136").Value = -1165765.5
